$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.884.57"
$ws.Range("E2").Value = "  +3.13%  "
$ws.Range("D3").Value = "3.794.47"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "696.40"
$ws.Range("E5").Value = "  +10.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.35"
$ws.Range("E6").Value = "  +5.14%  "
$ws.Range("D7").Value = "3.792.09"
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("E10").Value = "  +3.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.44"
$ws.Range("E11").Value = "  +8.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").Value = "  +8.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.28"
$ws.Range("E14").Value = "  +4.58%  "
$ws.Range("D15").Value = "4.435.73"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "3.798.35"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "70.941.91"
$ws.Range("E17").Value = "  +3.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.83"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.20"
$ws.Range("E19").Value = "  +3.36%  "
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.09"
$ws.Range("E21").Value = "  +17.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "484.50"
$ws.Range("E22").Value = "  +3.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.713"
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.19"
$ws.Range("E24").Value = "  +3.30%  "
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.42"
$ws.Range("E26").Value = "  +3.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.45"
$ws.Range("E27").Value = "  +4.21%  "
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("D29").Value = "3.947.37"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.06"
$ws.Range("E31").Value = "  +15.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.28"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.51"
$ws.Range("E33").Value = "  +5.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.60"
$ws.Range("E34").Value = "  +4.68%  "
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.22"
$ws.Range("E36").Value = "  +4.28%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "3.747.14"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("E39").Value = "  +2.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.50"
$ws.Range("E40").Value = "  +8.80%  "
$ws.Range("E41").Value = "  +4.30%  "
$ws.Range("E42").Value = "  +14.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.971"
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("E44").Value = "  +23.23%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "163.06"
$ws.Range("E47").Value = "  +4.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.27"
$ws.Range("E48").Value = "  +4.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.78"
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.300"
$ws.Range("E50").Value = "  +2.82%  "
$ws.Range("E51").Value = "  -1.01%  "
